# storeitems.xlsx - "fixed look and feel, updated user balance from db in all forms"
#
# The "store" worksheet is an append-only purchase log (A=item name,
# B=price, C=timestamp, D=user id, E=picture path). This adds the five new
# purchase records that were logged after the last export, picking up
# right after the existing row 23.
#
# Note: some of the "user id" values (e.g. "1234", "123123123") look
# numeric but must be stored as text, matching every other row in column D
# (and matching how the sheet already stores them). A leading apostrophe
# forces text entry the same way a user typing into Excel would, and the
# Style reset keeps the cell on the workbook's single default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("store")

function Set-TextCell($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 24
$ws.Range("A24").Value = "Buzz Big Doll"
$ws.Range("B24").Value = 125
Set-TextCell $ws.Range("C24") "08/08/2024 21:16:48"
Set-TextCell $ws.Range("D24") "1234"
$ws.Range("E24").Value = "pics/buzz_toy.png"

# Row 25
$ws.Range("A25").Value = "Mickeymouse Sticker"
$ws.Range("B25").Value = 20
Set-TextCell $ws.Range("C25") "08/08/2024 21:16:52"
Set-TextCell $ws.Range("D25") "1234"
$ws.Range("E25").Value = "pics/mickeymousesticker.png"

# Row 26
$ws.Range("A26").Value = "Ironman Sticker"
$ws.Range("B26").Value = 20
Set-TextCell $ws.Range("C26") "10/08/2024 19:49:28"
Set-TextCell $ws.Range("D26") "123123123"
$ws.Range("E26").Value = "pics/ironmansticker.png"

# Row 27
$ws.Range("A27").Value = "Big Hello Kitty Sticker"
$ws.Range("B27").Value = 40
Set-TextCell $ws.Range("C27") "10/08/2024 20:14:47"
Set-TextCell $ws.Range("D27") "123123123"
$ws.Range("E27").Value = "pics/hellokittysticker.png"

# Row 28
$ws.Range("A28").Value = "Truck Toy"
$ws.Range("B28").Value = 80
Set-TextCell $ws.Range("C28") "10/08/2024 20:14:50"
Set-TextCell $ws.Range("D28") "123123123"
$ws.Range("E28").Value = "pics/truck_toy.png"
